# Apply the "10.2.1" workbook update:
#  - Add the 2021-column (J) figures for rows 19-26 (numeric, formatted as "0.0",
#    matching the style already used by column I in those rows).
#  - Fill in the previously-empty J27 total cell.
#  - Move the active selection to N8 (as last left by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J values for rows 19-26.
$newValues = @{
    19 = 12.434613462352335
    20 = 16.80050595536094
    21 = 11.282963378125267
    22 = 25.042808754677555
    23 = 3.2011163356916352
    24 = 13.523574517571838
    25 = 6.1196997869329204
    26 = 5.9488136666578013
}

foreach ($row in $newValues.Keys) {
    $cell = $ws.Cells.Item($row, 10)   # column J
    $cell.Value = $newValues[$row]
    $cell.NumberFormat = "0.0"
}

# J27 cell already exists (bottom total row) but was empty - just set its value.
$ws.Range("J27").Value = 5.2451982064110645

# Update the sheet's active cell/selection.
$ws.Range("N8").Select() | Out-Null
